$d = $word.ActiveDocument

# 1. "Week 7 Questions (26 pts)" -> "... (24 pts)"
$d.Content.Find.Execute("(26 pts)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(24 pts)", 2)

# 2. Move the "_GoBack" bookmark from after "is a" (".docx" extension sentence)
#    to inside "be shared." (between "be sha" and "red.")
$r = $d.Content
$r.Find.Execute("be sha", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$pt = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $pt)

# 3. "Provide a brief overview ... (10 pts)" -> "... (8 pts)"
$d.Content.Find.Execute("(10 pts)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(8 pts)", 2)
